$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (e.g. '217.14')
# must be forced to Text format first, otherwise Excel auto-converts them
# to floating point numbers (losing exact formatting / trailing zeros).

$ws.Range('D2').Value = '27.166.70'
$ws.Range('E2').Value = '  +0.01%  '

$ws.Range('D3').Value = '1.638.47'
$ws.Range('E3').Value = '  -0.55%  '

$ws.Range('E4').Value = '  +0.28%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '217.14'
$ws.Range('E5').Value = '  -0.71%  '

$ws.Range('E6').Value = '  +1.59%  '

$ws.Range('E7').Value = '  +0.22%  '

$ws.Range('E8').Value = '  -0.74%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0626'
$ws.Range('E9').Value = '  -0.17%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.12'
$ws.Range('E10').Value = '  -1.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0849'
$ws.Range('E11').Value = '  +0.20%  '

$ws.Range('E12').Value = '  -0.48%  '

$ws.Range('D13').Value = '1.623.31'
$ws.Range('E13').Value = '  -1.53%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.14'
$ws.Range('E14').Value = '  +0.15%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.545'
$ws.Range('E15').Value = '  +0.58%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.34'
$ws.Range('E16').Value = '  -1.95%  '

$ws.Range('D17').Value = '27.166.35'
$ws.Range('E17').Value = '  +0.11%  '

$ws.Range('E18').Value = '  +0.20%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '216.82'
$ws.Range('E19').Value = '  -2.56%  '

$ws.Range('E20').Value = '  +0.27%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.88'
$ws.Range('E21').Value = '  +1.22%  '

$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.43'
$ws.Range('E22').Value = '  -0.56%  '

$ws.Range('B23').Value = 'Toncoin'
$ws.Range('C23').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '2.53'
$ws.Range('E23').Value = '  +4.04%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.15'
$ws.Range('E24').Value = '  -1.29%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '147.76'
$ws.Range('E25').Value = '  +0.19%  '

$ws.Range('E26').Value = '  +0.32%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.39'
$ws.Range('E27').Value = '  -0.35%  '

$ws.Range('E28').Value = '  -0.65%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.68'
$ws.Range('E29').Value = '  -1.11%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0509'
$ws.Range('E30').Value = '  +0.32%  '

$ws.Range('E31').Value = '  -0.71%  '

$ws.Range('E32').Value = '  +1.19%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.03'
$ws.Range('E33').Value = '  -0.46%  '

$ws.Range('D34').Value = '1.304.44'
$ws.Range('E34').Value = '  +2.38%  '

$ws.Range('E35').Value = '  -0.29%  '

$ws.Range('E36').Value = '  +1.01%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.549'
$ws.Range('E38').Value = '  +0.74%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.855'
$ws.Range('E39').Value = '  +0.88%  '

$ws.Range('E40').Value = '  +0.22%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.27'
$ws.Range('E41').Value = '  +5.41%  '

$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('E43').Value = '  -0.58%  '

$ws.Range('D44').Value = '1.777.77'
$ws.Range('E44').Value = '  -0.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '62.62'
$ws.Range('E45').Value = '  -1.72%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.87'
$ws.Range('E46').Value = '  -1.98%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('E47').Value = '  -1.08%  '

$ws.Range('E48').Value = '  +2.33%  '

$ws.Range('E49').Value = '  +0.14%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.63'
$ws.Range('E50').Value = '  -1.65%  '

$ws.Range('B51').Value = 'Algorand'
$ws.Range('C51').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0960'
$ws.Range('E51').Value = '  -1.44%  '
